$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report week dates) ---
$ws.Range("A8").Characters(21,2).Text = "28"
$ws.Range("C9").Characters(27,8).Text = "7/10/2023"
$ws.Range("C9").Characters(47,8).Text = "7/16/2023"

# --- Weekly crime-stat table updates (rows 14-29) ---
# Row 14
$ws.Range("D14").Copy($ws.Range("C14"))
# Row 15
$ws.Range("F14").Copy($ws.Range("C15"))
$ws.Range("D14").Copy($ws.Range("D15"))
$ws.Range("E14").Copy($ws.Range("E15"))
$ws.Range("F15").Value = 3
$ws.Range("H15").Value = -40
$ws.Range("I15").Value = 11
$ws.Range("K15").Value = -42.105263157894
$ws.Range("L15").Value = -26.666666666666
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = -31.25
# Row 16
$ws.Range("C16").Value = 10
$ws.Range("D16").Value = 11
$ws.Range("E16").Value = -9.090909090909
$ws.Range("F16").Value = 36
$ws.Range("G16").Value = 33
$ws.Range("H16").Value = 9.090909090909
$ws.Range("I16").Value = 168
$ws.Range("J16").Value = 145
$ws.Range("K16").Value = 15.862068965517
$ws.Range("L16").Value = 30.232558139534
$ws.Range("M16").Value = -5.084745762711
$ws.Range("N16").Value = -75.722543352601
# Row 17
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 80
$ws.Range("F17").Value = 35
$ws.Range("G17").Value = 27
$ws.Range("H17").Value = 29.629629629629
$ws.Range("I17").Value = 250
$ws.Range("J17").Value = 247
$ws.Range("K17").Value = 1.214574898785
$ws.Range("L17").Value = 31.578947368421
$ws.Range("M17").Value = 46.198830409356
$ws.Range("N17").Value = 10.132158590308
# Row 18
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 19
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = 90
$ws.Range("I18").Value = 88
$ws.Range("J18").Value = 80
$ws.Range("K18").Value = 10
$ws.Range("L18").Value = 27.536231884058
$ws.Range("M18").Value = -46.341463414634
$ws.Range("N18").Value = -92.177777777777
# Row 19
$ws.Range("C19").Value = 22
$ws.Range("D19").Value = 20
$ws.Range("E19").Value = 10
$ws.Range("F19").Value = 92
$ws.Range("G19").Value = 63
$ws.Range("H19").Value = 46.031746031746
$ws.Range("I19").Value = 450
$ws.Range("J19").Value = 518
$ws.Range("K19").Value = -13.127413127413
$ws.Range("L19").Value = 66.051660516605
$ws.Range("M19").Value = 82.926829268292
$ws.Range("N19").Value = -43.75
# Row 20
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 16
$ws.Range("E20").Value = -56.25
$ws.Range("F20").Value = 27
$ws.Range("G20").Value = 44
$ws.Range("H20").Value = -38.636363636363
$ws.Range("I20").Value = 165
$ws.Range("J20").Value = 161
$ws.Range("K20").Value = 2.484472049689
$ws.Range("L20").Value = 52.777777777777
$ws.Range("M20").Value = 30.952380952381
$ws.Range("N20").Value = -86.307053941908
# Row 21
$ws.Range("C21").Value = 53
$ws.Range("D21").Value = 55
$ws.Range("E21").Value = -3.636363636363
$ws.Range("F21").Value = 214
$ws.Range("G21").Value = 183
$ws.Range("H21").Value = 16.939890710382
$ws.Range("I21").Value = 1134
$ws.Range("J21").Value = 1173
$ws.Range("K21").Value = -3.324808184143
$ws.Range("L21").Value = 44.458598726114
$ws.Range("M21").Value = 26.5625
$ws.Range("N21").Value = -72.171779141104
# Row 22
$ws.Range("F14").Copy($ws.Range("C22"))
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 3
$ws.Range("E22").Value = -66.666666666666
$ws.Range("G22").Value = 8
$ws.Range("H22").Value = -62.5
$ws.Range("I22").Value = 56
$ws.Range("J22").Value = 28
$ws.Range("K22").Value = 100
$ws.Range("L22").Value = 700
$ws.Range("M22").Value = 211.111111111111
# Row 24
$ws.Range("D24").Value = 10
$ws.Range("E24").Value = 240
$ws.Range("F24").Value = 142
$ws.Range("G24").Value = 138
$ws.Range("H24").Value = 2.898550724637
$ws.Range("I24").Value = 1097
$ws.Range("J24").Value = 964
$ws.Range("K24").Value = 13.796680497925
$ws.Range("L24").Value = 49.658935879945
$ws.Range("M24").Value = 90.121317157712
# Row 25
$ws.Range("C25").Value = 21
$ws.Range("D25").Value = 13
$ws.Range("E25").Value = 61.538461538461
$ws.Range("F25").Value = 83
$ws.Range("H25").Value = 12.162162162162
$ws.Range("I25").Value = 503
$ws.Range("J25").Value = 490
$ws.Range("K25").Value = 2.653061224489
$ws.Range("L25").Value = 18.632075471698
$ws.Range("M25").Value = 1.411290322580
# Row 26
$ws.Range("C26").Value = 2
$ws.Range("D14").Copy($ws.Range("D26"))
$ws.Range("E14").Copy($ws.Range("E26"))
$ws.Range("F26").Value = 4
$ws.Range("H26").Value = -20
$ws.Range("I26").Value = 25
$ws.Range("K26").Value = -7.407407407407
$ws.Range("L26").Value = -7.407407407407
# Row 27
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = 66.666666666666
$ws.Range("F27").Value = 15
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 150
$ws.Range("I27").Value = 84
$ws.Range("J27").Value = 52
$ws.Range("K27").Value = 61.538461538461
$ws.Range("L27").Value = 25.373134328358
# Row 28
$ws.Range("D14").Copy($ws.Range("C28"))
$ws.Range("N28").Value = -88.571428571428
# Row 29
$ws.Range("D14").Copy($ws.Range("C29"))
$ws.Range("N29").Value = -90.625
